$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A and set a date/time number format style for it
# (target stored width 14.85546875 chars; engine rounds ColumnWidth to 1/6
# character steps, so 14 is the closest achievable input)
$ws.Columns.Item(1).ColumnWidth = 14

# Header cell A1 should use the same style (numFmtId 22) as the new date cell below.
# We set the number format on the whole column A data cells that will carry dates.
$ws.Range("A1").NumberFormat = "m/d/yy h:mm"

# Add the new row of data
$ws.Range("A2").Value = 42605.6475462963
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 46
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "Bag"
